$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('AA2').Value = 'maa://21246 (91.2), maa://36684 (98.59), ***maa://22731 (6.67)'
$ws.Range('AE2').Value = 'maa://25251 (92.5), ***maa://21730 (17.19), ***maa://39501 (25.0), *maa://36675 (60.0)'
$ws.Range('G3').Value = 'maa://21247 (98.28), *maa://22748 (75.0)'
$ws.Range('K3').Value = '*maa://22880 (69.93), maa://20276 (82.86), *maa://22749 (62.5)'
$ws.Range('W3').Value = 'maa://27396 (84.91), maa://27484 (95.74), maa://27480 (82.35)'
$ws.Range('AE3').Value = '*maa://21289 (70.0)'
$ws.Range('C4').Value = 'maa://24632 (93.23), **maa://24303 (36.36), maa://22499 (85.71), maa://22746 (100.0)'
$ws.Range('W4').Value = '**maa://32495 (47.54), ***maa://31785 (18.02), ***maa://36683 (26.67)'
$ws.Range('C5').Value = 'maa://21245 (82.29), maa://22744 (83.33)'
$ws.Range('C6').Value = '*maa://42407 (75.0)'
$ws.Range('K7').Value = 'maa://28624 (91.55), maa://24957 (97.44)'
$ws.Range('O7').Value = 'maa://22750 (97.06)'
$ws.Range('S7').Value = 'maa://21291 (89.47)'
$ws.Range('W7').Value = 'maa://22399 (94.62), *maa://22758 (71.15)'
$ws.Range('O8').Value = 'maa://32931 (88.61), *maa://21916 (60.34), maa://23252 (92.31), **maa://22759 (45.45), maa://37496 (100.0)'
$ws.Range('W8').Value = 'maa://21411 (96.01)'
$ws.Range('AE8').Value = '*maa://24479 (76.71), *maa://21990 (53.85)'
$ws.Range('O9').Value = 'maa://22736 (80.25)'
$ws.Range('W9').Value = 'maa://26223 (96.94)'
$ws.Range('AE9').Value = 'maa://26206 (89.16), **maa://22865 (45.65)'
$ws.Range('C10').Value = '***maa://25695 (19.19), **maa://32237 (37.84), ***maa://34206 (18.18), ***maa://39951 (19.23), **maa://39243 (33.33)'
$ws.Range('O10').Value = 'maa://28977 (94.74), *maa://23264 (62.96), maa://36669 (87.5)'
$ws.Range('S10').Value = 'maa://27395 (96.03), maa://22755 (87.62), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('W10').Value = 'maa://22301 (97.41), maa://22726 (100.0)'
$ws.Range('AE10').Value = '*maa://25021 (56.94), *maa://22733 (58.62), maa://22761 (100.0)'
$ws.Range('C11').Value = 'maa://36707 (99.64)'
$ws.Range('AE11').Value = 'maa://31203 (94.44), ***maa://24394 (19.23)'
$ws.Range('G12').Value = 'maa://21867 (90.0)'
$ws.Range('W12').Value = 'maa://22753 (91.84), *maa://21485 (76.74), maa://37962 (81.25)'
$ws.Range('AA12').Value = 'maa://23669 (95.51), maa://36677 (92.5), maa://39872 (84.62)'
$ws.Range('C13').Value = 'maa://24999 (91.45), maa://36673 (91.8), maa://25001 (85.51)'
$ws.Range('G13').Value = '*maa://21248 (75.73), **maa://22728 (47.62)'
$ws.Range('W13').Value = '*maa://34957 (75.56), *maa://22768 (53.33)'
$ws.Range('AE13').Value = '**maa://22737 (30.6), maa://39883 (87.5), *maa://39885 (73.68)'
$ws.Range('C14').Value = 'maa://30764 (86.36)'
$ws.Range('K14').Value = 'maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (93.02)'
$ws.Range('O14').Value = 'maa://23250 (98.48), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('R14').Value = 1
$ws.Range('S14').Value = 'maa://22521 (94.44)'
$ws.Range('W14').Value = 'maa://37468 (92.86)'
$ws.Range('C15').Value = '*maa://22743 (76.19), maa://22734 (83.33), *maa://30808 (64.29), ***maa://36048 (13.33)'
$ws.Range('G15').Value = 'maa://24304 (88.46), maa://21478 (91.18)'
$ws.Range('O15').Value = 'maa://24762 (89.58), *maa://22727 (70.0)'
$ws.Range('AE15').Value = 'maa://21364 (80.61), *maa://22766 (73.0), *maa://36666 (78.46)'
$ws.Range('C16').Value = 'maa://21441 (96.17), maa://36679 (93.94), maa://37650 (95.45)'
$ws.Range('W16').Value = 'maa://28501 (97.44), maa://28051 (95.83)'
$ws.Range('AE16').Value = '*maa://23911 (61.96), maa://27755 (91.78)'
$ws.Range('C18').Value = 'maa://24570 (96.55)'
$ws.Range('G18').Value = 'maa://24421 (90.57)'
$ws.Range('K18').Value = 'maa://22466 (88.55), *maa://22732 (51.85)'
$ws.Range('W18').Value = 'maa://21917 (97.5), maa://22741 (83.33)'
$ws.Range('AE18').Value = '*maa://24313 (57.62), **maa://29784 (46.15)'
$ws.Range('S19').Value = 'maa://24386 (98.8)'
$ws.Range('AA19').Value = '*maa://30709 (60.8), *maa://36668 (52.17)'
$ws.Range('K20').Value = 'maa://41331 (88.89)'
$ws.Range('S20').Value = 'maa://29113 (95.45)'
$ws.Range('AA21').Value = '*maa://21443 (78.72), **maa://23820 (30.91)'
$ws.Range('W22').Value = 'maa://21282 (98.83), *maa://37649 (66.67)'
$ws.Range('C23').Value = '***maa://28036 (28.36), *maa://41753 (75.0)'
$ws.Range('K23').Value = 'maa://39756 (92.36), maa://39875 (95.65)'
$ws.Range('S23').Value = 'maa://24387 (82.86), maa://31212 (95.83)'
$ws.Range('C24').Value = 'maa://24368 (80.36)'
$ws.Range('W24').Value = 'maa://23504 (92.9), maa://29988 (86.21), **maa://22892 (40.14), *maa://25141 (76.86), maa://36663 (80.36), ***maa://22815 (23.08)'
$ws.Range('AE24').Value = 'maa://22523 (85.03), *maa://36672 (76.74), maa://29910 (94.0), **maa://21440 (34.55)'
$ws.Range('C25').Value = 'maa://29753 (95.13)'
$ws.Range('AA26').Value = '*maa://42235 (64.29)'
$ws.Range('G27').Value = '**maa://21283 (48.65), maa://34494 (100.0), **maa://36665 (44.44), *maa://39601 (80.0)'
$ws.Range('C28').Value = 'maa://24465 (90.35), maa://25725 (82.28)'
$ws.Range('K28').Value = '*maa://30770 (79.07)'
$ws.Range('W28').Value = 'maa://39929 (86.15), ***maa://39723 (14.71), maa://41749 (81.25)'
$ws.Range('AE28').Value = 'maa://36660 (93.77), *maa://36701 (64.0)'
$ws.Range('AE29').Value = '*maa://24080 (68.96), ***maa://34960 (9.09)'
$ws.Range('K30').Value = 'maa://30442 (94.44)'
$ws.Range('G32').Value = 'maa://21895 (97.01), maa://36667 (98.08), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('S32').Value = 'maa://41108 (90.91), maa://41238 (94.44)'
$ws.Range('AD32').Value = 1
$ws.Range('AE32').Value = 'maa://42408 (100.0)'
$ws.Range('K35').Value = 'maa://41296 (97.96)'
$ws.Range('G39').Value = 'maa://25199 (86.11), maa://36670 (88.06), maa://30434 (87.27), ***maa://25036 (16.0)'
$ws.Range('O39').Value = 'maa://24709 (92.0)'
$ws.Range('G41').Value = 'maa://24466 (95.0)'
$ws.Range('O41').Value = '**maa://35616 (37.93)'
$ws.Range('G43').Value = 'maa://22525 (92.62), maa://21284 (82.93)'
$ws.Range('S44').Value = 'maa://39366 (86.96)'
$ws.Range('G46').Value = 'maa://35931 (92.51)'
$ws.Range('G51').Value = '*maa://30769 (80.0)'
$ws.Range('G58').Value = '*maa://37964 (63.16)'
